$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1252.9688
$ws.Range("I33").Value = 273.13635
$ws.Range("J33").Value = 3408.6
$ws.Range("K33").Value = 273.13635
$ws.Range("L33").Value = 3408.6
$ws.Range("M33").Value = -44.13634999999999
$ws.Range("N33").Value = -3866.6

$ws.Range("H70").Value = 1124.5
$ws.Range("I70").Value = 998
$ws.Range("J70").Value = 1166.6666
$ws.Range("K70").Value = 2994
$ws.Range("L70").Value = 3499.9998
$ws.Range("M70").Value = -2724
$ws.Range("N70").Value = -4039.9998

$ws.Range("H73").Value = 1124.5
$ws.Range("I73").Value = 998
$ws.Range("J73").Value = 1166.6666
$ws.Range("K73").Value = 2994
$ws.Range("L73").Value = 3499.9998
$ws.Range("M73").Value = -2058
$ws.Range("N73").Value = -5371.9998

$ws.Range("H88").Value = 2220.5557
$ws.Range("I88").Value = 618.1667
$ws.Range("J88").Value = 2678.3809
$ws.Range("K88").Value = 618.1667
$ws.Range("L88").Value = 2678.3809
$ws.Range("M88").Value = -212.1667
$ws.Range("N88").Value = -3490.3809

$ws.Range("H91").Value = 2220.5557
$ws.Range("I91").Value = 618.1667
$ws.Range("J91").Value = 2678.3809
$ws.Range("K91").Value = 618.1667
$ws.Range("L91").Value = 2678.3809
$ws.Range("M91").Value = 785.8333
$ws.Range("N91").Value = -5486.3809

$ws.Range("H94").Value = 83334320
$ws.Range("I94").Value = 1077.7273
$ws.Range("K94").Value = 1077.7273
$ws.Range("M94").Value = -626.7273

$ws.Range("H111").Value = 10113.9375
$ws.Range("I111").Value = 18482.428
$ws.Range("J111").Value = 3605.111
$ws.Range("K111").Value = 55447.284
$ws.Range("L111").Value = 10815.333
$ws.Range("M111").Value = -52380.284
$ws.Range("N111").Value = -16949.333

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H133").Value = 38398
$ws.Range("J133").Value = 38398
$ws.Range("L133").Value = 38398
$ws.Range("N133").Value = -48518

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35967.707
$ws.Range("I32").Value = 6496.073
$ws.Range("K32").Value = 6496.073
$ws.Range("M32").Value = -6209.073

$ws.Range("H122").Value = 2564.1428
$ws.Range("I122").Value = 2044
$ws.Range("J122").Value = 2954.25
$ws.Range("K122").Value = 6132
$ws.Range("L122").Value = 8862.75
$ws.Range("M122").Value = -3682
$ws.Range("N122").Value = -13762.75

$ws.Range("H139").Value = 48810.57
$ws.Range("I139").Value = 30000
$ws.Range("J139").Value = 51945.668
$ws.Range("K139").Value = 30000
$ws.Range("L139").Value = 51945.668
$ws.Range("M139").Value = -24860
$ws.Range("N139").Value = -62225.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1583.1666
$ws.Range("I99").Value = 1317.6666
$ws.Range("K99").Value = 1317.6666
$ws.Range("M99").Value = 180.3334

$ws.Range("H134").Value = 2459.5
$ws.Range("I134").Value = 2596.6667
$ws.Range("K134").Value = 7790.000100000001
$ws.Range("M134").Value = -5255.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31443.572
$ws.Range("I31").Value = 78548.69500000001
$ws.Range("J31").Value = 3608.7273
$ws.Range("K31").Value = 78548.69500000001
$ws.Range("L31").Value = 3608.7273
$ws.Range("M31").Value = -78253.69500000001
$ws.Range("N31").Value = -4198.7273

$ws.Range("H34").Value = 31443.572
$ws.Range("I34").Value = 78548.69500000001
$ws.Range("J34").Value = 3608.7273
$ws.Range("K34").Value = 78548.69500000001
$ws.Range("L34").Value = 3608.7273
$ws.Range("M34").Value = -78346.69500000001
$ws.Range("N34").Value = -4012.7273

$ws.Range("H94").Value = 1780
$ws.Range("J94").Value = 1650
$ws.Range("L94").Value = 1650
$ws.Range("N94").Value = -2552

$ws.Range("H122").Value = 874.2857
$ws.Range("I122").Value = 858.4
$ws.Range("J122").Value = 914
$ws.Range("K122").Value = 2575.2
$ws.Range("L122").Value = 2742
$ws.Range("M122").Value = -125.1999999999998
$ws.Range("N122").Value = -7642

$ws.Range("H134").Value = 1339.2646
$ws.Range("I134").Value = 1319.2963
$ws.Range("J134").Value = 1416.2858
$ws.Range("K134").Value = 3957.8889
$ws.Range("L134").Value = 4248.857400000001
$ws.Range("M134").Value = -1422.8889
$ws.Range("N134").Value = -9318.857400000001

$ws.Range("H141").Value = 61616.07
$ws.Range("I141").Value = 150000
$ws.Range("J141").Value = 37511.363
$ws.Range("K141").Value = 150000
$ws.Range("L141").Value = 37511.363
$ws.Range("M141").Value = -144820
$ws.Range("N141").Value = -47871.363

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 24555.25
$ws.Range("J37").Value = 24555.25
$ws.Range("L37").Value = 73665.75
$ws.Range("N37").Value = -73889.75

$ws.Range("H133").Value = 11600
$ws.Range("J133").Value = 11600
$ws.Range("L133").Value = 34800
$ws.Range("N133").Value = -44920

$ws.Range("H141").Value = 2851.3076
$ws.Range("I141").Value = 2373.4
$ws.Range("J141").Value = 4444.3335
$ws.Range("K141").Value = 7120.200000000001
$ws.Range("L141").Value = 13333.0005
$ws.Range("M141").Value = -1940.200000000001
$ws.Range("N141").Value = -23693.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 700
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2100
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 350
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1851.8889
$ws.Range("J93").Value = 1114
$ws.Range("L93").Value = 1114
$ws.Range("N93").Value = -3610

$ws.Range("H122").Value = 5091.8
$ws.Range("I122").Value = 4966.6665
$ws.Range("J122").Value = 5279.5
$ws.Range("K122").Value = 14899.9995
$ws.Range("L122").Value = 15838.5
$ws.Range("M122").Value = -12449.9995
$ws.Range("N122").Value = -20738.5

$ws.Range("H132").Value = 3575.6956
$ws.Range("I132").Value = 3611.85
$ws.Range("J132").Value = 3334.6667
$ws.Range("K132").Value = 10835.55
$ws.Range("L132").Value = 10004.0001
$ws.Range("M132").Value = -8305.549999999999
$ws.Range("N132").Value = -15064.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 250001600
$ws.Range("I96").Value = 500002000
$ws.Range("J96").Value = 1199
$ws.Range("K96").Value = 500002000
$ws.Range("L96").Value = 1199
$ws.Range("M96").Value = -500000627
$ws.Range("N96").Value = -3945

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 1889.5
$ws.Range("I122").Value = 1299
$ws.Range("J122").Value = 2243.8
$ws.Range("K122").Value = 3897
$ws.Range("L122").Value = 6731.400000000001
$ws.Range("M122").Value = -1447
$ws.Range("N122").Value = -11631.4

$ws.Range("H132").Value = 761.1724
$ws.Range("I132").Value = 558.14813
$ws.Range("J132").Value = 3502
$ws.Range("K132").Value = 1674.44439
$ws.Range("L132").Value = 10506
$ws.Range("M132").Value = 855.5556099999999
$ws.Range("N132").Value = -15566

$ws.Range("H135").Value = 39428.145
$ws.Range("J135").Value = 39428.145
$ws.Range("L135").Value = 39428.145
$ws.Range("N135").Value = -49568.145
